$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right after row 223 (date 44917,
# "Superior Seedless"). Every existing record from the old row 224 onward
# shifts down by one row (old 224 -> new 225, ..., old 235 -> new 236).
$ws.Rows.Item(224).EntireRow.Insert()

# Populate the newly inserted row 224 with the new record's data.
$ws.Cells.Item(224, 1).Value  = 2
$ws.Cells.Item(224, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(224, 3).Value  = "Coquimbo"
$ws.Cells.Item(224, 4).Value  = 44917
$ws.Cells.Item(224, 5).Value  = 4
$ws.Cells.Item(224, 6).Value  = "Fruta"
$ws.Cells.Item(224, 7).Value  = 100109
$ws.Cells.Item(224, 8).Value  = "Uva"
$ws.Cells.Item(224, 9).Value  = 100109001
$ws.Cells.Item(224, 10).Value = "Uva"
$ws.Cells.Item(224, 11).Value = "Superior Seedless"
$ws.Cells.Item(224, 12).Value = "Primera"
$ws.Cells.Item(224, 13).Value = 340
$ws.Cells.Item(224, 14).Value = 13500
$ws.Cells.Item(224, 15).Value = 14000
$ws.Cells.Item(224, 16).Value = 13750
$ws.Cells.Item(224, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(224, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(224, 19).Value = 1375
$ws.Cells.Item(224, 20).Value = 10
